# Weekly update: new Fukumoto price records for Vega Modelo de Temuco - Naranja.
# Insert 3 new rows at the top of the existing price block (row 744), shifting
# the existing rows 744-786 down to 747-789, then populate the 3 new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new blank rows above the current row 744 (shifts 744:786 -> 747:789)
$ws.Rows("744:746").Insert()

# Row 744
$ws.Range("A744").Value = 10
$ws.Range("B744").Value = "Vega Modelo de Temuco"
$ws.Range("C744").Value = "La Araucanía"
$ws.Range("D744").Value = 44706
$ws.Range("E744").Value = 9
$ws.Range("F744").Value = "Fruta"
$ws.Range("G744").Value = 100102
$ws.Range("H744").Value = "Cítricos"
$ws.Range("I744").Value = 100102005
$ws.Range("J744").Value = "Naranja"
$ws.Range("K744").Value = "Fukumoto"
$ws.Range("L744").Value = "Primera"
$ws.Range("M744").Value = 5
$ws.Range("N744").Value = 250000
$ws.Range("O744").Value = 250000
$ws.Range("P744").Value = 250000
$ws.Range("Q744").Value = "$/bins (400 kilos)"
$ws.Range("R744").Value = "Región de O'Higgins"
$ws.Range("S744").Value = 625
$ws.Range("T744").Value = 400

# Row 745
$ws.Range("A745").Value = 10
$ws.Range("B745").Value = "Vega Modelo de Temuco"
$ws.Range("C745").Value = "La Araucanía"
$ws.Range("D745").Value = 44706
$ws.Range("E745").Value = 9
$ws.Range("F745").Value = "Fruta"
$ws.Range("G745").Value = 100102
$ws.Range("H745").Value = "Cítricos"
$ws.Range("I745").Value = 100102005
$ws.Range("J745").Value = "Naranja"
$ws.Range("K745").Value = "Fukumoto"
$ws.Range("L745").Value = "Segunda"
$ws.Range("M745").Value = 8
$ws.Range("N745").Value = 200000
$ws.Range("O745").Value = 200000
$ws.Range("P745").Value = 200000
$ws.Range("Q745").Value = "$/bins (400 kilos)"
$ws.Range("R745").Value = "Región de O'Higgins"
$ws.Range("S745").Value = 500
$ws.Range("T745").Value = 400

# Row 746
$ws.Range("A746").Value = 10
$ws.Range("B746").Value = "Vega Modelo de Temuco"
$ws.Range("C746").Value = "La Araucanía"
$ws.Range("D746").Value = 44706
$ws.Range("E746").Value = 9
$ws.Range("F746").Value = "Fruta"
$ws.Range("G746").Value = 100102
$ws.Range("H746").Value = "Cítricos"
$ws.Range("I746").Value = 100102005
$ws.Range("J746").Value = "Naranja"
$ws.Range("K746").Value = "Fukumoto"
$ws.Range("L746").Value = "Segunda"
$ws.Range("M746").Value = 250
$ws.Range("N746").Value = 6000
$ws.Range("O746").Value = 6000
$ws.Range("P746").Value = 6000
$ws.Range("Q746").Value = "$/malla 15 kilos"
$ws.Range("R746").Value = "Región de O'Higgins"
$ws.Range("S746").Value = 400
$ws.Range("T746").Value = 15
